# The "Rede" column (D) values "Total (4)" and "Pública (4)" had the
# trailing " (4)" footnote-reference marker removed, becoming plain
# "Total" and "Pública". Scan every used row in column D and strip the
# suffix wherever it is found (rows already reading "Total"/"Pública",
# e.g. the tail of the sheet, are left untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($val -eq "Total (4)") {
        $cell.Value = "Total"
    } elseif ($val -eq "Pública (4)") {
        $cell.Value = "Pública"
    }
}
